# Insert a new column before column B ("Colecao"), splitting the old
# combined "ID" text column (e.g. "MTR 12450") into a collection-prefix
# column (B) and a purely numeric ID column (C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh column at B - shifts old B (ID text) -> C, C->D, D->E, etc.
$ws.Columns("B").Insert()

# Row data: Colecao prefix (col B) and numeric ID (col C, replacing the
# old combined "PREFIX NNNN" text that lived there before the insert).
$colecao = @("MTR", "MTR", "UFMG", "UFMG", "LACV", "LACV", "CHUNB", "CHUNB", "CHUNB", "CHUNB", "CHUNB")
$ids = @(12450, 16143, 1327, 2175, 3288, 3236, 67393, 832, 63198, 11514, 37088)

for ($i = 0; $i -lt $colecao.Length; $i++) {
    $row = $i + 2
    $bCell = $ws.Cells.Item($row, 2)
    $bCell.Value = $colecao[$i]
    $bCell.HorizontalAlignment = -4108   # xlCenter - matches column B/C's shared style

    $ws.Cells.Item($row, 3).Value = $ids[$i]
}

# New header for the inserted column, added last (B1 already inherited
# style s="2" from the insert, matching the other header cells).
$ws.Range("B1").Value = "Colecao"

# The leftover " " placeholder used to live at B17; after the column
# insert it shifted to C17. It has no place in the new layout - remove it
# (and its formatting) entirely.
$ws.Range("C17").Clear()

# Match the author's final selection.
$ws.Range("F12").Select()
